$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet "Sheet3". The target state adds a
# new "Sheet1" worksheet in front of it (so tab order becomes Sheet1, Sheet3)
# and "Sheet1" becomes the active/selected sheet.
#
# Worksheets.Add() inserts the new sheet before the currently active sheet
# and activates it, which is exactly the placement/activation we need here.
$sheet1 = $wb.Worksheets.Add()
$sheet1.Name = "Sheet1"

# Populate the new Sheet1 with its data/formulas:
#   A1=100  B1=1            C1==12+12 (24)
#   A2=200  B2=2
#   A3=300  B3=3
#   A4==SUM(A1:A3) (600)    B4==SUM(B1:B3) (6)
$sheet1.Range("A1").Value = 100
$sheet1.Range("B1").Value = 1
$sheet1.Range("C1").Formula = "=12+12"

$sheet1.Range("A2").Value = 200
$sheet1.Range("B2").Value = 2

$sheet1.Range("A3").Value = 300
$sheet1.Range("B3").Value = 3

$sheet1.Range("A4").Formula = "=SUM(A1:A3)"
$sheet1.Range("B4").Formula = "=SUM(B1:B3)"

# Match the recorded selection on the new active sheet (cell C2 selected).
$sheet1.Range("C2").Select()
